$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 366.34375
$ws.Range("I19").Value = 354.0909
$ws.Range("M19").Value = -179.0909
$ws.Range("K19").Value = 354.0909
$ws.Range("H28").Value = 6622.15
$ws.Range("L28").Value = 2524
$ws.Range("I28").Value = 9354.25
$ws.Range("N28").Value = -3494
$ws.Range("J28").Value = 2524
$ws.Range("M28").Value = -8869.25
$ws.Range("K28").Value = 9354.25
$ws.Range("J40").Value = 4554.7
$ws.Range("L40").Value = 4554.7
$ws.Range("H40").Value = 3268.1304
$ws.Range("N40").Value = -4904.7
$ws.Range("I40").Value = 2278.4614
$ws.Range("M40").Value = -2103.4614
$ws.Range("K40").Value = 2278.4614
$ws.Range("K51").Value = 6996.3335
$ws.Range("H51").Value = 6909
$ws.Range("J51").Value = 6843.5
$ws.Range("M51").Value = -6512.3335
$ws.Range("I51").Value = 6996.3335
$ws.Range("N51").Value = -7811.5
$ws.Range("L51").Value = 6843.5
$ws.Range("J61").Value = 0
$ws.Range("N61").Value = ""
$ws.Range("L61").Value = 0
$ws.Range("H61").Value = 1179.125
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("I92").Value = 1002.36365
$ws.Range("H92").Value = 1002.36365
$ws.Range("M92").Value = 245.63635
$ws.Range("K92").Value = 1002.36365
$ws.Range("J92").Value = 0
$ws.Range("H98").Value = 33221.11
$ws.Range("I98").Value = 40127.645
$ws.Range("J98").Value = 9048.25
$ws.Range("K98").Value = 40127.645
$ws.Range("N98").Value = -12044.25
$ws.Range("M98").Value = -38629.645
$ws.Range("L98").Value = 9048.25
$ws.Range("N100").Value = -82748
$ws.Range("H100").Value = 22862982
$ws.Range("M100").Value = -50200019
$ws.Range("K100").Value = 50200560
$ws.Range("J100").Value = 81666
$ws.Range("L100").Value = 81666
$ws.Range("I100").Value = 50200560
$ws.Range("K106").Value = 15436543
$ws.Range("H106").Value = 9502095
$ws.Range("I106").Value = 15436543
$ws.Range("M106").Value = -15435912
$ws.Range("H112").Value = 3186.6
$ws.Range("L112").Value = 11122.35
$ws.Range("J112").Value = 3707.45
$ws.Range("N112").Value = -13338.35
$ws.Range("H122").Value = 33221.11
$ws.Range("L122").Value = 27144.75
$ws.Range("K122").Value = 120382.935
$ws.Range("N122").Value = -32044.75
$ws.Range("M122").Value = -117932.935
$ws.Range("I122").Value = 40127.645
$ws.Range("J122").Value = 9048.25
$ws.Range("M127").Value = 1533.1426
$ws.Range("K127").Value = 3426.8574
$ws.Range("H127").Value = 1142.2858
$ws.Range("I127").Value = 1142.2858
$ws.Range("K132").Value = 11691
$ws.Range("I132").Value = 3897
$ws.Range("M132").Value = -9161
$ws.Range("H132").Value = 3931.0962
$ws.Range("H137").Value = 7599.729
$ws.Range("N137").Value = -17594.7861
$ws.Range("J137").Value = 4164.9287
$ws.Range("L137").Value = 12494.7861
$ws.Range("I141").Value = 4559.857
$ws.Range("M141").Value = -8499.571
$ws.Range("K141").Value = 13679.571
$ws.Range("H141").Value = 4747.448

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L30").Value = 2000
$ws.Range("N30").Value = -2300
$ws.Range("H30").Value = 2000
$ws.Range("J30").Value = 2000
$ws.Range("I32").Value = 5444.65
$ws.Range("H32").Value = 5520.512
$ws.Range("M32").Value = -5157.65
$ws.Range("K32").Value = 5444.65
$ws.Range("K45").Value = 340950
$ws.Range("I45").Value = 340950
$ws.Range("N45").Value = -12087.333
$ws.Range("H45").Value = 231077.78
$ws.Range("L45").Value = 11333.333
$ws.Range("M45").Value = -340573
$ws.Range("J45").Value = 11333.333
$ws.Range("H46").Value = 9504.666999999999
$ws.Range("L46").Value = 8405.6
$ws.Range("N46").Value = -9043.6
$ws.Range("J46").Value = 8405.6
$ws.Range("J61").Value = 4219.2
$ws.Range("N61").Value = -4643.2
$ws.Range("I61").Value = 3483.2683
$ws.Range("K61").Value = 3483.2683
$ws.Range("L61").Value = 4219.2
$ws.Range("M61").Value = -3271.2683
$ws.Range("H61").Value = 3563.261
$ws.Range("I64").Value = 45000
$ws.Range("K64").Value = 45000
$ws.Range("H64").Value = 45000
$ws.Range("M64").Value = -44752
$ws.Range("M67").Value = -44142
$ws.Range("I67").Value = 45000
$ws.Range("K67").Value = 45000
$ws.Range("H67").Value = 45000
$ws.Range("H74").Value = 4033.05
$ws.Range("J74").Value = 18333.666
$ws.Range("I74").Value = 1509.4117
$ws.Range("K74").Value = 1509.4117
$ws.Range("M74").Value = -635.4117000000001
$ws.Range("N74").Value = -20081.666
$ws.Range("L74").Value = 18333.666
$ws.Range("K77").Value = 7547.058500000001
$ws.Range("N77").Value = -100404.33
$ws.Range("M77").Value = -3179.058500000001
$ws.Range("I77").Value = 1509.4117
$ws.Range("J77").Value = 18333.666
$ws.Range("H77").Value = 4033.05
$ws.Range("L77").Value = 91668.33
$ws.Range("M102").Value = -13227
$ws.Range("H102").Value = 10648.305
$ws.Range("I102").Value = 14849
$ws.Range("K102").Value = 14849
$ws.Range("H107").Value = 15000
$ws.Range("J107").Value = 15000
$ws.Range("L107").Value = 15000
$ws.Range("N107").Value = -22680
$ws.Range("K110").Value = 684.9231
$ws.Range("H110").Value = 1544.5555
$ws.Range("I110").Value = 684.9231
$ws.Range("M110").Value = 1360.0769
$ws.Range("H122").Value = 1077254
$ws.Range("L122").Value = 11018727
$ws.Range("K122").Value = 9569.793600000001
$ws.Range("N122").Value = -11023627
$ws.Range("M122").Value = -7119.793600000001
$ws.Range("I122").Value = 3189.9312
$ws.Range("J122").Value = 3672909
$ws.Range("K132").Value = 23089.731
$ws.Range("I132").Value = 7696.577
$ws.Range("M132").Value = -20559.731
$ws.Range("H132").Value = 6908.9487
$ws.Range("M136").Value = -7899.804900000001
$ws.Range("K136").Value = 10449.8049
$ws.Range("J136").Value = 4219.2
$ws.Range("H136").Value = 3563.261
$ws.Range("L136").Value = 12657.6
$ws.Range("N136").Value = -17757.6
$ws.Range("I136").Value = 3483.2683

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = ""
$ws.Range("H57").Value = 49000
$ws.Range("J57").Value = 0
$ws.Range("N93").Value = -68244
$ws.Range("H93").Value = 64500
$ws.Range("J93").Value = 64500
$ws.Range("L93").Value = 64500
$ws.Range("I97").Value = 4875
$ws.Range("K97").Value = 4875
$ws.Range("H97").Value = 28250
$ws.Range("M97").Value = -3884
$ws.Range("J105").Value = 11333.333
$ws.Range("N105").Value = -14827.333
$ws.Range("H105").Value = 95555.45
$ws.Range("K105").Value = 127138.75
$ws.Range("I105").Value = 127138.75
$ws.Range("L105").Value = 11333.333
$ws.Range("M105").Value = -125391.75
$ws.Range("H107").Value = 2459.4167
$ws.Range("J107").Value = 3954
$ws.Range("L107").Value = 3954
$ws.Range("I107").Value = 2160.5
$ws.Range("M107").Value = -240.5
$ws.Range("N107").Value = -7794
$ws.Range("K107").Value = 2160.5
$ws.Range("N116").Value = -48928
$ws.Range("L116").Value = 39750
$ws.Range("H116").Value = 39750
$ws.Range("J116").Value = 39750
$ws.Range("H119").Value = 19975
$ws.Range("J119").Value = 19975
$ws.Range("N119").Value = -29651
$ws.Range("L119").Value = 19975
$ws.Range("K134").Value = 39576.714
$ws.Range("H134").Value = 12368.75
$ws.Range("I134").Value = 13192.238
$ws.Range("M134").Value = -37041.714
$ws.Range("J135").Value = 99999
$ws.Range("L135").Value = 99999
$ws.Range("N135").Value = -110139
$ws.Range("H135").Value = 99999
$ws.Range("J136").Value = 0
$ws.Range("H136").Value = 49000
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = ""
$ws.Range("H137").Value = 79999
$ws.Range("N137").Value = -90199
$ws.Range("J137").Value = 79999
$ws.Range("L137").Value = 79999
$ws.Range("J140").Value = 76732.25
$ws.Range("L140").Value = 76732.25
$ws.Range("N140").Value = -87092.25
$ws.Range("H140").Value = 76732.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M86").Value = -9477
$ws.Range("I86").Value = 10600
$ws.Range("H86").Value = 13746.615
$ws.Range("K86").Value = 10600
$ws.Range("L86").Value = 14318.728
$ws.Range("N86").Value = -16564.728
$ws.Range("J86").Value = 14318.728
$ws.Range("J89").Value = 14318.728
$ws.Range("K89").Value = 53000
$ws.Range("H89").Value = 13746.615
$ws.Range("L89").Value = 71593.64
$ws.Range("I89").Value = 10600
$ws.Range("M89").Value = -47384
$ws.Range("N89").Value = -82825.64
$ws.Range("H94").Value = 2628.2666
$ws.Range("N94").Value = -3927.2
$ws.Range("K94").Value = 1834.4
$ws.Range("J94").Value = 3025.2
$ws.Range("I94").Value = 1834.4
$ws.Range("M94").Value = -1383.4
$ws.Range("L94").Value = 3025.2
$ws.Range("N99").Value = -8232.684000000001
$ws.Range("H99").Value = 147292.14
$ws.Range("L99").Value = 5236.684
$ws.Range("M99").Value = -314485
$ws.Range("K99").Value = 315983
$ws.Range("I99").Value = 315983
$ws.Range("J99").Value = 5236.684
$ws.Range("H107").Value = 8318.379000000001
$ws.Range("J107").Value = 229
$ws.Range("L107").Value = 229
$ws.Range("I107").Value = 9251.77
$ws.Range("M107").Value = -7331.77
$ws.Range("N107").Value = -4069
$ws.Range("K107").Value = 9251.77
$ws.Range("N109").Value = ""
$ws.Range("L109").Value = 0
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("H112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("N112").Value = ""
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = ""
$ws.Range("H122").Value = 8523.529
$ws.Range("L122").Value = 5518.2858
$ws.Range("K122").Value = 39607.2
$ws.Range("N122").Value = -10418.2858
$ws.Range("M122").Value = -37157.2
$ws.Range("I122").Value = 13202.4
$ws.Range("J122").Value = 1839.4286
$ws.Range("J126").Value = 5236.684
$ws.Range("I126").Value = 315983
$ws.Range("L126").Value = 15710.052
$ws.Range("K126").Value = 947949
$ws.Range("H126").Value = 147292.14
$ws.Range("M126").Value = -945479
$ws.Range("N126").Value = -20650.052
$ws.Range("K132").Value = 6413.4786
$ws.Range("I132").Value = 2137.8262
$ws.Range("M132").Value = -3883.4786
$ws.Range("H132").Value = 14649.25
$ws.Range("K134").Value = 9229.2855
$ws.Range("H134").Value = 3003.9412
$ws.Range("I134").Value = 3076.4285
$ws.Range("M134").Value = -6694.2855

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 22617880
$ws.Range("I4").Value = 573021.4399999999
$ws.Range("K4").Value = 1719064.32
$ws.Range("M4").Value = -1718952.32
$ws.Range("N5").Value = -2312496.2
$ws.Range("K5").Value = 655.875
$ws.Range("L5").Value = 2312272.2
$ws.Range("I5").Value = 218.625
$ws.Range("M5").Value = -543.875
$ws.Range("J5").Value = 770757.4
$ws.Range("H5").Value = 345632.56
$ws.Range("H12").Value = 143.64285
$ws.Range("K12").Value = 744.4285500000001
$ws.Range("N12").Value = -463.428571
$ws.Range("L12").Value = 117.428571
$ws.Range("I12").Value = 248.14285
$ws.Range("M12").Value = -571.4285500000001
$ws.Range("J12").Value = 39.142857
$ws.Range("K17").Value = 3483
$ws.Range("L17").Value = 6462.428400000001
$ws.Range("I17").Value = 1161
$ws.Range("M17").Value = -3314
$ws.Range("N17").Value = -6800.428400000001
$ws.Range("H17").Value = 1933.4445
$ws.Range("J17").Value = 2154.1428
$ws.Range("L32").Value = 2439
$ws.Range("J32").Value = 813
$ws.Range("H32").Value = 813
$ws.Range("N32").Value = -3005
$ws.Range("J44").Value = 4113.778
$ws.Range("I44").Value = 285
$ws.Range("L44").Value = 12341.334
$ws.Range("H44").Value = 3417.6365
$ws.Range("N44").Value = -13137.334
$ws.Range("K44").Value = 855
$ws.Range("M44").Value = -457
$ws.Range("H58").Value = 3455.4443
$ws.Range("N68").Value = -93779468
$ws.Range("H68").Value = 26323882
$ws.Range("K68").Value = 5248.9998
$ws.Range("J68").Value = 31259282
$ws.Range("I68").Value = 1749.6666
$ws.Range("L68").Value = 93777846
$ws.Range("M68").Value = -4437.9998
$ws.Range("K71").Value = 15746.9994
$ws.Range("I71").Value = 1749.6666
$ws.Range("M71").Value = -11690.9994
$ws.Range("J71").Value = 31259282
$ws.Range("L71").Value = 281333538
$ws.Range("N71").Value = -281341650
$ws.Range("H71").Value = 26323882
$ws.Range("K80").Value = 17982.4995
$ws.Range("L80").Value = 238566.81
$ws.Range("H80").Value = 53571.176
$ws.Range("N80").Value = -240438.81
$ws.Range("I80").Value = 5994.1665
$ws.Range("M80").Value = -17046.4995
$ws.Range("J80").Value = 79522.27
$ws.Range("L83").Value = 715700.4300000001
$ws.Range("N83").Value = -725060.4300000001
$ws.Range("I83").Value = 5994.1665
$ws.Range("J83").Value = 79522.27
$ws.Range("H83").Value = 53571.176
$ws.Range("M83").Value = -49267.4985
$ws.Range("K83").Value = 53947.4985
$ws.Range("K132").Value = 1800
$ws.Range("I132").Value = 200
$ws.Range("M132").Value = 730
$ws.Range("H132").Value = 100439
$ws.Range("J135").Value = 770757.4
$ws.Range("I135").Value = 218.625
$ws.Range("M135").Value = 567.375
$ws.Range("L135").Value = 6936816.600000001
$ws.Range("N135").Value = -6941886.600000001
$ws.Range("K135").Value = 1967.625
$ws.Range("H135").Value = 345632.56
$ws.Range("I139").Value = 6668368
$ws.Range("K139").Value = 20005104
$ws.Range("H139").Value = 4447012
$ws.Range("M139").Value = -19999964

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("N19").Value = ""
$ws.Range("M19").Value = ""
$ws.Range("K19").Value = 0
$ws.Range("H94").Value = 89814760
$ws.Range("N94").Value = -112019802
$ws.Range("J94").Value = 112018450
$ws.Range("L94").Value = 112018450
$ws.Range("I97").Value = 6404.5
$ws.Range("K97").Value = 6404.5
$ws.Range("H97").Value = 5391.037
$ws.Range("M97").Value = -5908.5
$ws.Range("M102").Value = -7573.352999999999
$ws.Range("H102").Value = 6409.9644
$ws.Range("I102").Value = 9195.352999999999
$ws.Range("K102").Value = 9195.352999999999
$ws.Range("J126").Value = 11489.223
$ws.Range("I126").Value = 11886.308
$ws.Range("L126").Value = 34467.669
$ws.Range("K126").Value = 35658.924
$ws.Range("H126").Value = 11723.863
$ws.Range("M126").Value = -33188.924
$ws.Range("N126").Value = -39407.669
$ws.Range("K132").Value = 4807.174199999999
$ws.Range("L132").Value = 5997
$ws.Range("I132").Value = 1602.3914
$ws.Range("J132").Value = 1999
$ws.Range("M132").Value = -2277.174199999999
$ws.Range("N132").Value = -11057
$ws.Range("H132").Value = 1673.2142
$ws.Range("N133").Value = -80115.664
$ws.Range("J133").Value = 69995.664
$ws.Range("H133").Value = 69995.664
$ws.Range("L133").Value = 69995.664

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K7").Value = 24978.13
$ws.Range("M7").Value = -24866.13
$ws.Range("H7").Value = 17372.756
$ws.Range("I7").Value = 24978.13
$ws.Range("L16").Value = 4197.5
$ws.Range("K16").Value = 5269.75
$ws.Range("J16").Value = 4197.5
$ws.Range("N16").Value = -4537.5
$ws.Range("I16").Value = 5269.75
$ws.Range("H16").Value = 5001.6875
$ws.Range("M16").Value = -5099.75
$ws.Range("I22").Value = 18841.545
$ws.Range("J22").Value = 1736.7
$ws.Range("K22").Value = 18841.545
$ws.Range("L22").Value = 1736.7
$ws.Range("M22").Value = -18546.545
$ws.Range("H22").Value = 10696.381
$ws.Range("N22").Value = -2326.7
$ws.Range("L27").Value = 1736.7
$ws.Range("M27").Value = -18734.545
$ws.Range("I27").Value = 18841.545
$ws.Range("N27").Value = -1950.7
$ws.Range("K27").Value = 18841.545
$ws.Range("H27").Value = 10696.381
$ws.Range("J27").Value = 1736.7
$ws.Range("I31").Value = 5012.5
$ws.Range("H31").Value = 3508.3333
$ws.Range("N31").Value = -996
$ws.Range("L31").Value = 500
$ws.Range("M31").Value = -4764.5
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 5012.5
$ws.Range("J40").Value = 12412.385
$ws.Range("L40").Value = 12412.385
$ws.Range("H40").Value = 22723.46
$ws.Range("N40").Value = -12684.385
$ws.Range("I40").Value = 33034.54
$ws.Range("M40").Value = -32898.54
$ws.Range("K40").Value = 33034.54
$ws.Range("N100").Value = -4154.25
$ws.Range("H100").Value = 2898.6667
$ws.Range("M100").Value = -2218.8
$ws.Range("K100").Value = 2759.8
$ws.Range("J100").Value = 3072.25
$ws.Range("L100").Value = 3072.25
$ws.Range("I100").Value = 2759.8
$ws.Range("I126").Value = 24978.13
$ws.Range("K126").Value = 74934.39
$ws.Range("H126").Value = 17372.756
$ws.Range("M126").Value = -72464.39
$ws.Range("J129").Value = 68429
$ws.Range("L129").Value = 68429
$ws.Range("H129").Value = 68429
$ws.Range("N129").Value = -78429
$ws.Range("K132").Value = 3736003.2
$ws.Range("L132").Value = 19087.3638
$ws.Range("I132").Value = 1245334.4
$ws.Range("J132").Value = 6362.4546
$ws.Range("M132").Value = -3733473.2
$ws.Range("N132").Value = -24147.3638
$ws.Range("H132").Value = 652782.5600000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 10000000
$ws.Range("I26").Value = 10000000
$ws.Range("K26").Value = 10000000
$ws.Range("M26").Value = -9999707
$ws.Range("L62").Value = 5775
$ws.Range("H62").Value = 287732
$ws.Range("J62").Value = 5775
$ws.Range("N62").Value = -7023
$ws.Range("L65").Value = 28875
$ws.Range("J65").Value = 5775
$ws.Range("H65").Value = 287732
$ws.Range("N65").Value = -35115
$ws.Range("N99").Value = ""
$ws.Range("H99").Value = 23431.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -20436.5
$ws.Range("K99").Value = 23431.5
$ws.Range("I99").Value = 23431.5
$ws.Range("J99").Value = 0
$ws.Range("N100").Value = -175970
$ws.Range("M100").Value = -111458.2
$ws.Range("K100").Value = 111999.2
$ws.Range("J100").Value = 87444
$ws.Range("L100").Value = 174888
$ws.Range("I100").Value = 55999.6
$ws.Range("H122").Value = 4558.452
$ws.Range("K122").Value = 5923.0344
$ws.Range("M122").Value = -3473.0344
$ws.Range("I122").Value = 1974.3448
$ws.Range("I126").Value = 25359.059
$ws.Range("K126").Value = 76077.177
$ws.Range("H126").Value = 17666.629
$ws.Range("M126").Value = -73607.177
$ws.Range("K130").Value = 0
$ws.Range("H130").Value = 86143
$ws.Range("M130").Value = ""
$ws.Range("I130").Value = 0
$ws.Range("K132").Value = 46064.325
$ws.Range("L132").Value = 17181.5448
$ws.Range("I132").Value = 15354.775
$ws.Range("J132").Value = 5727.1816
$ws.Range("M132").Value = -43534.325
$ws.Range("N132").Value = -22241.5448
$ws.Range("H132").Value = 13278.235
$ws.Range("M136").Value = -944104.14
$ws.Range("K136").Value = 946654.14
$ws.Range("J136").Value = 3725.3572
$ws.Range("H136").Value = 246256.7
$ws.Range("L136").Value = 11176.0716
$ws.Range("N136").Value = -16276.0716
$ws.Range("I136").Value = 315551.38
